$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column F (dSF) values to match repulled data
$ws.Range("F2").Value = -6
$ws.Range("F3").Value = 8
$ws.Range("F4").Value = -4
$ws.Range("F6").Value = 1
$ws.Range("F9").Value = -8
$ws.Range("F10").Value = -3
